# Auto-generated script applying scheduled price-refresh updates
# to the Typhon_Profits workbook (per-sheet Leve profit tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 757.1667
$ws.Range("I28").Value = 303.5
$ws.Range("K28").Value = 303.5
$ws.Range("M28").Value = 181.5
$ws.Range("H62").Value = 5906.3076
$ws.Range("I62").Value = 4971.25
$ws.Range("J62").Value = 6321.8887
$ws.Range("K62").Value = 4971.25
$ws.Range("L62").Value = 6321.8887
$ws.Range("M62").Value = -4347.25
$ws.Range("N62").Value = -7569.8887
$ws.Range("H65").Value = 5906.3076
$ws.Range("I65").Value = 4971.25
$ws.Range("J65").Value = 6321.8887
$ws.Range("K65").Value = 24856.25
$ws.Range("L65").Value = 31609.4435
$ws.Range("M65").Value = -21736.25
$ws.Range("N65").Value = -37849.4435
$ws.Range("H76").Value = 2528300
$ws.Range("I76").Value = 3159
$ws.Range("J76").Value = 9262009
$ws.Range("K76").Value = 3159
$ws.Range("L76").Value = 9262009
$ws.Range("M76").Value = -2844
$ws.Range("N76").Value = -9262639
$ws.Range("H79").Value = 2528300
$ws.Range("I79").Value = 3159
$ws.Range("J79").Value = 9262009
$ws.Range("K79").Value = 3159
$ws.Range("L79").Value = 9262009
$ws.Range("M79").Value = -2067
$ws.Range("N79").Value = -9264193
$ws.Range("H80").Value = 11075891
$ws.Range("I80").Value = 432.33334
$ws.Range("J80").Value = 18743516
$ws.Range("K80").Value = 1297.00002
$ws.Range("L80").Value = 56230548
$ws.Range("M80").Value = -299.0000199999999
$ws.Range("N80").Value = -56232544
$ws.Range("H83").Value = 11075891
$ws.Range("I83").Value = 432.33334
$ws.Range("J83").Value = 18743516
$ws.Range("K83").Value = 3891.00006
$ws.Range("L83").Value = 168691644
$ws.Range("M83").Value = 1100.99994
$ws.Range("N83").Value = -168701628
$ws.Range("H86").Value = 34680.668
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 51021
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 51021
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -53267
$ws.Range("H89").Value = 34680.668
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 51021
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 255105
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -266337
$ws.Range("H94").Value = 2782.5
$ws.Range("I94").Value = 2782.5
$ws.Range("K94").Value = 2782.5
$ws.Range("M94").Value = -2331.5
$ws.Range("H132").Value = 3023.2632
$ws.Range("I132").Value = 3036.1333
$ws.Range("J132").Value = 2975
$ws.Range("K132").Value = 9108.3999
$ws.Range("L132").Value = 8925
$ws.Range("M132").Value = -6578.3999
$ws.Range("N132").Value = -13985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 14153.125
$ws.Range("I132").Value = 1287.8276
$ws.Range("J132").Value = 48070.727
$ws.Range("K132").Value = 3863.4828
$ws.Range("L132").Value = 144212.181
$ws.Range("M132").Value = -1333.4828
$ws.Range("N132").Value = -149272.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 802.5263
$ws.Range("I64").Value = 1013.2727
$ws.Range("J64").Value = 512.75
$ws.Range("K64").Value = 1013.2727
$ws.Range("L64").Value = 512.75
$ws.Range("M64").Value = -788.2727
$ws.Range("N64").Value = -962.75
$ws.Range("H67").Value = 802.5263
$ws.Range("I67").Value = 1013.2727
$ws.Range("J67").Value = 512.75
$ws.Range("K67").Value = 1013.2727
$ws.Range("L67").Value = 512.75
$ws.Range("M67").Value = -233.2727
$ws.Range("N67").Value = -2072.75
$ws.Range("H99").Value = 1152.3334
$ws.Range("I99").Value = 1188.1818
$ws.Range("J99").Value = 1096
$ws.Range("K99").Value = 1188.1818
$ws.Range("L99").Value = 1096
$ws.Range("M99").Value = 309.8181999999999
$ws.Range("N99").Value = -4092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 45916.668
$ws.Range("J68").Value = 45916.668
$ws.Range("L68").Value = 45916.668
$ws.Range("N68").Value = -47414.668
$ws.Range("H71").Value = 45916.668
$ws.Range("J71").Value = 45916.668
$ws.Range("L71").Value = 137750.004
$ws.Range("N71").Value = -145238.004
$ws.Range("H122").Value = 1475
$ws.Range("I122").Value = 1212.5
$ws.Range("K122").Value = 3637.5
$ws.Range("M122").Value = -1187.5
$ws.Range("H132").Value = 3307.36
$ws.Range("I132").Value = 2406
$ws.Range("K132").Value = 7218
$ws.Range("M132").Value = -4688
$ws.Range("H134").Value = 1580
$ws.Range("I134").Value = 950
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2850
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -315
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 862.875
$ws.Range("J98").Value = 600
$ws.Range("L98").Value = 1800
$ws.Range("N98").Value = -4796
$ws.Range("H131").Value = 703.97
$ws.Range("J131").Value = 727.72046
$ws.Range("L131").Value = 2183.16138
$ws.Range("N131").Value = -12263.16138

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 18034.098
$ws.Range("I132").Value = 1595.1052
$ws.Range("J132").Value = 44062.5
$ws.Range("K132").Value = 4785.3156
$ws.Range("L132").Value = 132187.5
$ws.Range("M132").Value = -2255.3156
$ws.Range("N132").Value = -137247.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 18000
$ws.Range("J51").Value = 18000
$ws.Range("L51").Value = 18000
$ws.Range("M51").Value = -18956
$ws.Range("H61").Value = 4229.6875
$ws.Range("I61").Value = 1436.5385
$ws.Range("J61").Value = 16333.333
$ws.Range("K61").Value = 1436.5385
$ws.Range("L61").Value = 16333.333
$ws.Range("M61").Value = -1234.5385
$ws.Range("N61").Value = -16737.333
$ws.Range("H113").Value = 4229.6875
$ws.Range("I113").Value = 1436.5385
$ws.Range("J113").Value = 16333.333
$ws.Range("K113").Value = 1436.5385
$ws.Range("L113").Value = 16333.333
$ws.Range("M113").Value = 733.4614999999999
$ws.Range("N113").Value = -20673.333
$ws.Range("H122").Value = 1156372.8
$ws.Range("I122").Value = 1636469.9
$ws.Range("K122").Value = 4909409.699999999
$ws.Range("M122").Value = -4906959.699999999
$ws.Range("H129").Value = 30000
$ws.Range("J129").Value = 30000
$ws.Range("L129").Value = 30000
$ws.Range("N129").Value = -40000
$ws.Range("H136").Value = 1640.9474
$ws.Range("I136").Value = 1441.2858
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 4323.857400000001
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -1773.857400000001
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 71429624
$ws.Range("I81").Value = 1193.4445
$ws.Range("J81").Value = 200000800
$ws.Range("K81").Value = 2386.889
$ws.Range("L81").Value = 400001600
$ws.Range("M81").Value = -1325.889
$ws.Range("N81").Value = -400003722
$ws.Range("H84").Value = 71429624
$ws.Range("I84").Value = 1193.4445
$ws.Range("J84").Value = 200000800
$ws.Range("K84").Value = 11934.445
$ws.Range("L84").Value = 2000008000
$ws.Range("M84").Value = -6630.445
$ws.Range("N84").Value = -2000018608
$ws.Range("H122").Value = 1246.8823
$ws.Range("I122").Value = 1257
$ws.Range("K122").Value = 3771
$ws.Range("M122").Value = -1321
